$d = $word.ActiveDocument
$d.Content.Find.Execute("39+37=76", $true, $false, $false, $false, $false, $true, 1, $false, "88-63=25", 2) | Out-Null
$d.Content.Find.Execute("46+3=49", $true, $false, $false, $false, $false, $true, 1, $false, "16+15=31", 2) | Out-Null
$d.Content.Find.Execute("4-2=2", $true, $false, $false, $false, $false, $true, 1, $false, "43-2=41", 2) | Out-Null
$d.Content.Find.Execute("93-63=30", $true, $false, $false, $false, $false, $true, 1, $false, "72+9=81", 2) | Out-Null
$d.Content.Find.Execute("15-11=4", $true, $false, $false, $false, $false, $true, 1, $false, "97-72=25", 2) | Out-Null
$d.Content.Find.Execute("65-49=16", $true, $false, $false, $false, $false, $true, 1, $false, "97-30=67", 2) | Out-Null
$d.Content.Find.Execute("36-15=21", $true, $false, $false, $false, $false, $true, 1, $false, "96-13=83", 2) | Out-Null
$d.Content.Find.Execute("73-50=23", $true, $false, $false, $false, $false, $true, 1, $false, "29+64=93", 2) | Out-Null
$d.Content.Find.Execute("67+3=70", $true, $false, $false, $false, $false, $true, 1, $false, "5-1=4", 2) | Out-Null
$d.Content.Find.Execute("20+71=91", $true, $false, $false, $false, $false, $true, 1, $false, "78-21=57", 2) | Out-Null
$d.Content.Find.Execute("31+23=54", $true, $false, $false, $false, $false, $true, 1, $false, "90-86=4", 2) | Out-Null
$d.Content.Find.Execute("4+1=5", $true, $false, $false, $false, $false, $true, 1, $false, "92-67=25", 2) | Out-Null
$d.Content.Find.Execute("41+5=46", $true, $false, $false, $false, $false, $true, 1, $false, "76+17=93", 2) | Out-Null
$d.Content.Find.Execute("80-44=36", $true, $false, $false, $false, $false, $true, 1, $false, "49+18=67", 2) | Out-Null
$d.Content.Find.Execute("15+63=78", $true, $false, $false, $false, $false, $true, 1, $false, "68+28=96", 2) | Out-Null
$d.Content.Find.Execute("64-19=45", $true, $false, $false, $false, $false, $true, 1, $false, "68-22=46", 2) | Out-Null
$d.Content.Find.Execute("49-3=46", $true, $false, $false, $false, $false, $true, 1, $false, "25+20=45", 2) | Out-Null
$d.Content.Find.Execute("71+10=81", $true, $false, $false, $false, $false, $true, 1, $false, "12+5=17", 2) | Out-Null
$d.Content.Find.Execute("63-0=63", $true, $false, $false, $false, $false, $true, 1, $false, "9+80=89", 2) | Out-Null
$d.Content.Find.Execute("54+10=64", $true, $false, $false, $false, $false, $true, 1, $false, "70+21=91", 2) | Out-Null
$d.Content.Find.Execute("41+42=83", $true, $false, $false, $false, $false, $true, 1, $false, "70-15=55", 2) | Out-Null
$d.Content.Find.Execute("89-52=37", $true, $false, $false, $false, $false, $true, 1, $false, "89-73=16", 2) | Out-Null
$d.Content.Find.Execute("90-1=89", $true, $false, $false, $false, $false, $true, 1, $false, "79-50=29", 2) | Out-Null
$d.Content.Find.Execute("57+36=93", $true, $false, $false, $false, $false, $true, 1, $false, "28+63=91", 2) | Out-Null
$d.Content.Find.Execute("34+44=78", $true, $false, $false, $false, $false, $true, 1, $false, "39+23=62", 2) | Out-Null
$d.Content.Find.Execute("51+23=74", $true, $false, $false, $false, $false, $true, 1, $false, "46-45=1", 2) | Out-Null
$d.Content.Find.Execute("34-21=13", $true, $false, $false, $false, $false, $true, 1, $false, "25+27=52", 2) | Out-Null
$d.Content.Find.Execute("0+16=16", $true, $false, $false, $false, $false, $true, 1, $false, "84-37=47", 2) | Out-Null
$d.Content.Find.Execute("24+41=65", $true, $false, $false, $false, $false, $true, 1, $false, "68+12=80", 2) | Out-Null
$d.Content.Find.Execute("83-47=36", $true, $false, $false, $false, $false, $true, 1, $false, "30-26=4", 2) | Out-Null
$d.Content.Find.Execute("75+5=80", $true, $false, $false, $false, $false, $true, 1, $false, "47+52=99", 2) | Out-Null
$d.Content.Find.Execute("74-24=50", $true, $false, $false, $false, $false, $true, 1, $false, "80-38=42", 2) | Out-Null
$d.Content.Find.Execute("92-57=35", $true, $false, $false, $false, $false, $true, 1, $false, "71-70=1", 2) | Out-Null
$d.Content.Find.Execute("83-67=16", $true, $false, $false, $false, $false, $true, 1, $false, "84-5=79", 2) | Out-Null
$d.Content.Find.Execute("28-3=25", $true, $false, $false, $false, $false, $true, 1, $false, "27+5=32", 2) | Out-Null
$d.Content.Find.Execute("98-48=50", $true, $false, $false, $false, $false, $true, 1, $false, "28+14=42", 2) | Out-Null
$d.Content.Find.Execute("44-11=33", $true, $false, $false, $false, $false, $true, 1, $false, "1+71=72", 2) | Out-Null
$d.Content.Find.Execute("84+5=89", $true, $false, $false, $false, $false, $true, 1, $false, "73+4=77", 2) | Out-Null
$d.Content.Find.Execute("18-5=13", $true, $false, $false, $false, $false, $true, 1, $false, "93-25=68", 2) | Out-Null
$d.Content.Find.Execute("5+14=19", $true, $false, $false, $false, $false, $true, 1, $false, "86-41=45", 2) | Out-Null
$d.Content.Find.Execute("37-3=34", $true, $false, $false, $false, $false, $true, 1, $false, "66-28=38", 2) | Out-Null
$d.Content.Find.Execute("22+45=67", $true, $false, $false, $false, $false, $true, 1, $false, "72-24=48", 2) | Out-Null
$d.Content.Find.Execute("44+9=53", $true, $false, $false, $false, $false, $true, 1, $false, "94-40=54", 2) | Out-Null
$d.Content.Find.Execute("54+6=60", $true, $false, $false, $false, $false, $true, 1, $false, "41-1=40", 2) | Out-Null
$d.Content.Find.Execute("53+31=84", $true, $false, $false, $false, $false, $true, 1, $false, "10+29=39", 2) | Out-Null
$d.Content.Find.Execute("71+27=98", $true, $false, $false, $false, $false, $true, 1, $false, "94-84=10", 2) | Out-Null
$d.Content.Find.Execute("56+23=79", $true, $false, $false, $false, $false, $true, 1, $false, "76-45=31", 2) | Out-Null
$d.Content.Find.Execute("15+84=99", $true, $false, $false, $false, $false, $true, 1, $false, "87-75=12", 2) | Out-Null
$d.Content.Find.Execute("63-38=25", $true, $false, $false, $false, $false, $true, 1, $false, "58-54=4", 2) | Out-Null
$d.Content.Find.Execute("71-1=70", $true, $false, $false, $false, $false, $true, 1, $false, "77-76=1", 2) | Out-Null
$d.Content.Find.Execute("13+67=80", $true, $false, $false, $false, $false, $true, 1, $false, "61+26=87", 2) | Out-Null
$d.Content.Find.Execute("63-57=6", $true, $false, $false, $false, $false, $true, 1, $false, "62+7=69", 2) | Out-Null
$d.Content.Find.Execute("21+61=82", $true, $false, $false, $false, $false, $true, 1, $false, "31+41=72", 2) | Out-Null
$d.Content.Find.Execute("60-37=23", $true, $false, $false, $false, $false, $true, 1, $false, "2+47=49", 2) | Out-Null
$d.Content.Find.Execute("6+70=76", $true, $false, $false, $false, $false, $true, 1, $false, "61-14=47", 2) | Out-Null
$d.Content.Find.Execute("63+35=98", $true, $false, $false, $false, $false, $true, 1, $false, "21+65=86", 2) | Out-Null
$d.Content.Find.Execute("39+27=66", $true, $false, $false, $false, $false, $true, 1, $false, "53-35=18", 2) | Out-Null
$d.Content.Find.Execute("20-14=6", $true, $false, $false, $false, $false, $true, 1, $false, "36-10=26", 2) | Out-Null
$d.Content.Find.Execute("71-20=51", $true, $false, $false, $false, $false, $true, 1, $false, "1+63=64", 2) | Out-Null
$d.Content.Find.Execute("34+45=79", $true, $false, $false, $false, $false, $true, 1, $false, "0+9=9", 2) | Out-Null
$d.Content.Find.Execute("48+0=48", $true, $false, $false, $false, $false, $true, 1, $false, "23-22=1", 2) | Out-Null
$d.Content.Find.Execute("21+62=83", $true, $false, $false, $false, $false, $true, 1, $false, "66-6=60", 2) | Out-Null
$d.Content.Find.Execute("88+6=94", $true, $false, $false, $false, $false, $true, 1, $false, "6+91=97", 2) | Out-Null
$d.Content.Find.Execute("0+54=54", $true, $false, $false, $false, $false, $true, 1, $false, "64-4=60", 2) | Out-Null
$d.Content.Find.Execute("1+11=12", $true, $false, $false, $false, $false, $true, 1, $false, "43-23=20", 2) | Out-Null
$d.Content.Find.Execute("65-26=39", $true, $false, $false, $false, $false, $true, 1, $false, "40+44=84", 2) | Out-Null
$d.Content.Find.Execute("57-47=10", $true, $false, $false, $false, $false, $true, 1, $false, "63-17=46", 2) | Out-Null
$d.Content.Find.Execute("67-15=52", $true, $false, $false, $false, $false, $true, 1, $false, "36+24=60", 2) | Out-Null
$d.Content.Find.Execute("88-59=29", $true, $false, $false, $false, $false, $true, 1, $false, "35+18=53", 2) | Out-Null
$d.Content.Find.Execute("16+31=47", $true, $false, $false, $false, $false, $true, 1, $false, "7+51=58", 2) | Out-Null
$d.Content.Find.Execute("10+71=81", $true, $false, $false, $false, $false, $true, 1, $false, "12+8=20", 2) | Out-Null
$d.Content.Find.Execute("71-53=18", $true, $false, $false, $false, $false, $true, 1, $false, "98-97=1", 2) | Out-Null
$d.Content.Find.Execute("99-42=57", $true, $false, $false, $false, $false, $true, 1, $false, "32+21=53", 2) | Out-Null
$d.Content.Find.Execute("89+2=91", $true, $false, $false, $false, $false, $true, 1, $false, "79-26=53", 2) | Out-Null
$d.Content.Find.Execute("14+2=16", $true, $false, $false, $false, $false, $true, 1, $false, "43-34=9", 2) | Out-Null
$d.Content.Find.Execute("80-64=16", $true, $false, $false, $false, $false, $true, 1, $false, "99-1=98", 2) | Out-Null
$d.Content.Find.Execute("68-7=61", $true, $false, $false, $false, $false, $true, 1, $false, "94-48=46", 2) | Out-Null
$d.Content.Find.Execute("56-47=9", $true, $false, $false, $false, $false, $true, 1, $false, "19+70=89", 2) | Out-Null
$d.Content.Find.Execute("16+75=91", $true, $false, $false, $false, $false, $true, 1, $false, "39+60=99", 2) | Out-Null
$d.Content.Find.Execute("57-22=35", $true, $false, $false, $false, $false, $true, 1, $false, "88-58=30", 2) | Out-Null
$d.Content.Find.Execute("74+9=83", $true, $false, $false, $false, $false, $true, 1, $false, "36-20=16", 2) | Out-Null
$d.Content.Find.Execute("66+15=81", $true, $false, $false, $false, $false, $true, 1, $false, "53-19=34", 2) | Out-Null
$d.Content.Find.Execute("3+70=73", $true, $false, $false, $false, $false, $true, 1, $false, "48+7=55", 2) | Out-Null
$d.Content.Find.Execute("8+16=24", $true, $false, $false, $false, $false, $true, 1, $false, "69-32=37", 2) | Out-Null
$d.Content.Find.Execute("83-63=20", $true, $false, $false, $false, $false, $true, 1, $false, "98-0=98", 2) | Out-Null
$d.Content.Find.Execute("34+5=39", $true, $false, $false, $false, $false, $true, 1, $false, "76-30=46", 2) | Out-Null
$d.Content.Find.Execute("8+24=32", $true, $false, $false, $false, $false, $true, 1, $false, "76+13=89", 2) | Out-Null
$d.Content.Find.Execute("56-38=18", $true, $false, $false, $false, $false, $true, 1, $false, "3+8=11", 2) | Out-Null
$d.Content.Find.Execute("49+22=71", $true, $false, $false, $false, $false, $true, 1, $false, "12+25=37", 2) | Out-Null
$d.Content.Find.Execute("86-6=80", $true, $false, $false, $false, $false, $true, 1, $false, "43+3=46", 2) | Out-Null
$d.Content.Find.Execute("87-73=14", $true, $false, $false, $false, $false, $true, 1, $false, "39+56=95", 2) | Out-Null
$d.Content.Find.Execute("65-9=56", $true, $false, $false, $false, $false, $true, 1, $false, "2+78=80", 2) | Out-Null
$d.Content.Find.Execute("43-36=7", $true, $false, $false, $false, $false, $true, 1, $false, "48-2=46", 2) | Out-Null
$d.Content.Find.Execute("83-20=63", $true, $false, $false, $false, $false, $true, 1, $false, "45-14=31", 2) | Out-Null
$d.Content.Find.Execute("14+48=62", $true, $false, $false, $false, $false, $true, 1, $false, "21+76=97", 2) | Out-Null
$d.Content.Find.Execute("58-1=57", $true, $false, $false, $false, $false, $true, 1, $false, "86+3=89", 2) | Out-Null
$d.Content.Find.Execute("35-28=7", $true, $false, $false, $false, $false, $true, 1, $false, "56+36=92", 2) | Out-Null
$d.Content.Find.Execute("40+2=42", $true, $false, $false, $false, $false, $true, 1, $false, "57-26=31", 2) | Out-Null
$d.Content.Find.Execute("99-38=61", $true, $false, $false, $false, $false, $true, 1, $false, "35+2=37", 2) | Out-Null
$d.Content.Find.Execute("88-1=87", $true, $false, $false, $false, $false, $true, 1, $false, "44+39=83", 2) | Out-Null
